$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the stepping table values in rows 61-63 (TOF PAC, TOF PAC +1kV, MCP)
# Columns C,D,E,F,G,I change from 16000 to 12000 (column H is already 12000)
foreach ($r in 61..63) {
    foreach ($col in @("C","D","E","F","G","I")) {
        $ws.Range("$col$r").Value = 12000
    }
}

# Update the active selection to reflect the last-edited cell
$null = $ws.Range("B15").Select()
